$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert two new rows before the totals row (row 20), shifting
#        the totals row (20->22) and the footer row (21->23) down.
$ws.Rows("20:21").Insert()

# Row heights for the two new rows (matches target layout)
$ws.Rows("20").RowHeight = 24.75
$ws.Rows("21").RowHeight = 25.5

# Clone formatting (styles + the merged-cell layout) of the last
# product row (19) into the two freshly inserted rows.
$ws.Range("A19:Q19").Copy()
$ws.Range("A20:Q20").PasteSpecial(-4122)
$ws.Range("A19:Q19").Copy()
$ws.Range("A21:Q21").PasteSpecial(-4122)

# Recreate the merged cells for the two new rows (Insert() does not
# carry the merge layout over automatically).
$ws.Range("A20:B20").Merge()
$ws.Range("C20:G20").Merge()
$ws.Range("H20:K20").Merge()
$ws.Range("L20:M20").Merge()
$ws.Range("N20:O20").Merge()

$ws.Range("A21:B21").Merge()
$ws.Range("C21:G21").Merge()
$ws.Range("H21:K21").Merge()
$ws.Range("L21:M21").Merge()
$ws.Range("N21:O21").Merge()

# --- 2. Re-write the product table rows 14-21 so that the two new
#        items ("SIMETHICONE-MUP ..." and "سهايه الجو") slot into
#        their correct alphabetical position, pushing the following
#        rows down by two. The leading "'" keeps these numeric-looking
#        strings stored as text (matching the source sheet's layout),
#        same as every other cell in this table.
$ws.Range("C14").Value = "SIMETHICONE-MUP 2% EMULSION ORAL DROPS 30 ML"
$ws.Range("H14").Value = "'2:0"
$ws.Range("L14").Value = "'1"
$ws.Range("N14").Value = "'27.00"
$ws.Range("P14").Value = "'27.0000"
$ws.Range("Q14").Value = "'1:0"

$ws.Range("C15").Value = "TERRAMYCIN EYE OINT. 5 GM"
$ws.Range("H15").Value = "'3:0"
$ws.Range("L15").Value = "'1"
$ws.Range("N15").Value = "'28.00"
$ws.Range("P15").Value = "'28.0000"
$ws.Range("Q15").Value = "'1:0"

$ws.Range("C16").Value = "TOBRIN 0.3% EYE DROPS 5 ML"
$ws.Range("H16").Value = "'2:0"
$ws.Range("L16").Value = "'1"
$ws.Range("N16").Value = "'23.00"
$ws.Range("P16").Value = "'23.0000"
$ws.Range("Q16").Value = "'1:0"

$ws.Range("C17").Value = "ايفا كيراتين حمام كريم"
$ws.Range("H17").Value = "'0:0"
$ws.Range("L17").Value = "'0"
$ws.Range("N17").Value = "'180.00"
$ws.Range("P17").Value = "'180.0000"
$ws.Range("Q17").Value = "'1:0"

$ws.Range("C18").Value = "سرنجات 3 سم"
$ws.Range("H18").Value = "'0:0"
$ws.Range("L18").Value = "'0"
$ws.Range("N18").Value = "'2.00"
$ws.Range("P18").Value = "'2.0000"
$ws.Range("Q18").Value = "'1:0"

$ws.Range("C19").Value = "سهايه الجو"
$ws.Range("H19").Value = "'1:0"
$ws.Range("L19").Value = "'0"
$ws.Range("N19").Value = "'15.00"
$ws.Range("P19").Value = "'15.0000"
$ws.Range("Q19").Value = "'1:0"

$ws.Range("A20").Value = 14
$ws.Range("C20").Value = "شامبو الفيف 400 مل"
$ws.Range("H20").Value = "'1:0"
$ws.Range("L20").Value = "'0"
$ws.Range("N20").Value = "'150.00"
$ws.Range("P20").Value = "'150.0000"
$ws.Range("Q20").Value = "'1:0"

$ws.Range("A21").Value = 15
$ws.Range("C21").Value = "ليفه"
$ws.Range("H21").Value = "'4:0"
$ws.Range("L21").Value = "'0"
$ws.Range("N21").Value = "'15.00"
$ws.Range("P21").Value = "'15.0000"
$ws.Range("Q21").Value = "'1:0"

# --- 3. Update the totals row (now row 22): total went from 691.43
#        to 733.43 (+27.00 +15.00 for the two new items).
$ws.Range("P22").Value = 733.42999999999995

# --- 4. Update the generated-at timestamp in the footer (now row 23).
$ws.Range("A23").Value = "Monday, 14 July, 2025 10:57 AM"
